# Adds a new "Database Create" test-case sheet right after the
# "Database GetConnection" sheet, mirroring its layout/header styling,
# and makes the new sheet the active tab (matching the authored edit).

$wb = $excel.ActiveWorkbook

# The sheet the new one is modeled on / inserted after.
$refSheet = $wb.Worksheets.Item("Database GetConnection")

# Remember where the selection was left on the reference sheet before
# the edit (Excel keeps per-sheet selection state even after the tab
# stops being active).
$refSheet.Range("A1:F1").Select() | Out-Null

# Insert the new worksheet right after "Database GetConnection".
$newWs = $wb.Worksheets.Add($null, $refSheet)
$newWs.Name = "Database Create"

# Copy the header row (values + formatting) from the reference sheet so
# the new sheet matches its look (fonts/borders/alignment/column widths).
$refSheet.Range("A1:F1").Copy($newWs.Range("A1:F1"))
$newWs.Rows(1).RowHeight = 27.6

# New test-case row for the "Database Create" method.
$newWs.Range("A2").Value = 1
$newWs.Range("B2").Value = "Valid"
$newWs.Range("C2").Value = "Creates Database"
$newWs.Range("D2").Value = "NA"
$newWs.Range("E2").Value = "No Exceptions (void)"

# Match column sizing used elsewhere in the workbook.
$newWs.Columns("C").ColumnWidth = 15.44140625
$newWs.Columns("E").ColumnWidth = 17.21875
$newWs.Columns("F").ColumnWidth = 10.33203125

# Leave the selection/active cell where the author last left it, and make
# the new sheet the active tab.
$newWs.Range("F2").Select() | Out-Null
$newWs.Activate()
